$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'insulated jogging pants for women',
    'insulated knee brace',
    'insulated leggings',
    'insulated leggings for women',
    'insulated leggings for women warm',
    'insulated leggings for women winter',
    'insulated leggins',
    'insulated riding pants',
    'insulated running gear',
    'insulated ski pants',
    'insulated skirt',
    'insulated skirt for women',
    'insulated skirts for women',
    'insulated track pants mens',
    'insulated white leggings',
    'insulated womens running pants',
    'insulated workout leggings for women',
    'insulator clothing',
    'insulator pants',
    'insulted pants for men',
    'insulted work pants',
    'intense moisture wrap',
    'intensity pants',
    'international clothing for women',
    'international womens day',
    'investments pants',
    'investments pants for women',
    'investments womens pants',
    'inward training',
    'iphone 5 skin wrap',
    'iphone 6 waistband for running',
    'iphone recovery',
    'iphone x s skin',
    'iphone x waistband for running',
    'iron neck pro',
    'it works body wrap tape',
    'itb band',
    'ivy park slide',
    'jacket and pants set for women nike',
    'january cross',
    'jersey legging',
    'jiggle balls for women',
    'jiggle key',
    'jiggle pants',
    'jockey capri',
    'jockey pants women',
    'jockey womens leggings',
    'jodhpurs women',
    'jogging tights women',
    'join pain',
    'join support',
    'joint 3x',
    'joint brace back brace',
    'joint braces',
    'joint comfort',
    'joint compression',
    'joint help',
    'joint muscle',
    'joint pain support',
    'joint performance',
    'joint power',
    'joint pro',
    'joint recovery',
    'joint rub',
    'joint stability',
    'joint supplements for runners',
    'joint support',
    'joint support and comfort',
    'joint support athlete',
    'joint support brace',
    'joint support compression pants',
    'joint support for athletes',
    'joint support for men',
    'joint support for runners',
    'joint support for women',
    'joint support knee brace',
    'joint support leggings',
    'joint support natural',
    'joint support now',
    'joint support plus',
    'joint support tape',
    'joint support women',
    'joint swelling',
    'joint tape',
    'joint works',
    'joint x',
    'june butter',
    'june leggings',
    'juniors pink leggings',
    'jusfitsu grey dress',
    'just live leggings',
    'just my size black pants',
    'just my size womens plus',
    'kentucky jodphurs',
    'key muscles of yoga',
    'keys art pro',
    'kilometers to miles',
    'kimchi blue dress',
    'kimchi blue dresses',
    'kimchi dress'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
